$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 189.33333
$ws.Range("I33").Value = 191.64706
$ws.Range("J33").Value = 150
$ws.Range("K33").Value = 191.64706
$ws.Range("L33").Value = 150
$ws.Range("M33").Value = 37.35293999999999
$ws.Range("N33").Value = -608

$ws.Range("H88").Value = 4216.75
$ws.Range("J88").Value = 4216.75
$ws.Range("L88").Value = 4216.75
$ws.Range("N88").Value = -5028.75

$ws.Range("H91").Value = 4216.75
$ws.Range("J91").Value = 4216.75
$ws.Range("L91").Value = 4216.75
$ws.Range("N91").Value = -7024.75

$ws.Range("H107").Value = 703.6774
$ws.Range("I107").Value = 364.5909
$ws.Range("J107").Value = 1532.5555
$ws.Range("K107").Value = 364.5909
$ws.Range("L107").Value = 1532.5555
$ws.Range("M107").Value = 1555.4091
$ws.Range("N107").Value = -5372.5555

$ws.Range("H132").Value = 1636.25
$ws.Range("I132").Value = 1364.7097
$ws.Range("J132").Value = 3319.8
$ws.Range("K132").Value = 4094.1291
$ws.Range("L132").Value = 9959.400000000001
$ws.Range("M132").Value = -1564.1291
$ws.Range("N132").Value = -15019.4

$ws.Range("H135").Value = 6313.5
$ws.Range("J135").Value = 16478.334
$ws.Range("L135").Value = 148305.006
$ws.Range("N135").Value = -153375.006

$ws.Range("H137").Value = 11906488
$ws.Range("I137").Value = 19231988
$ws.Range("J137").Value = 2552.125
$ws.Range("K137").Value = 57695964
$ws.Range("L137").Value = 7656.375
$ws.Range("M137").Value = -57693414
$ws.Range("N137").Value = -12756.375

$ws.Range("H138").Value = 2788.224
$ws.Range("J138").Value = 4824.3335
$ws.Range("L138").Value = 14473.0005
$ws.Range("N138").Value = -24753.0005

$ws.Range("H141").Value = 6953.5
$ws.Range("I141").Value = 6942.25
$ws.Range("K141").Value = 20826.75
$ws.Range("M141").Value = -15646.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 41665.668
$ws.Range("I31").Value = 12499.5
$ws.Range("J31").Value = 99998
$ws.Range("K31").Value = 12499.5
$ws.Range("L31").Value = 99998
$ws.Range("M31").Value = -12205.5
$ws.Range("N31").Value = -100586

$ws.Range("H74").Value = 1593.3636
$ws.Range("I74").Value = 1502.5714
$ws.Range("J74").Value = 3500
$ws.Range("K74").Value = 1502.5714
$ws.Range("L74").Value = 3500
$ws.Range("M74").Value = -628.5714
$ws.Range("N74").Value = -5248

$ws.Range("H77").Value = 1593.3636
$ws.Range("I77").Value = 1502.5714
$ws.Range("J77").Value = 3500
$ws.Range("K77").Value = 7512.857
$ws.Range("L77").Value = 17500
$ws.Range("M77").Value = -3144.857
$ws.Range("N77").Value = -26236

$ws.Range("H102").Value = 1065.48
$ws.Range("I102").Value = 1065.48
$ws.Range("K102").Value = 1065.48
$ws.Range("M102").Value = 556.52

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 1141.2
$ws.Range("I80").Value = 716
$ws.Range("J80").Value = 1295.8182
$ws.Range("K80").Value = 716
$ws.Range("L80").Value = 1295.8182
$ws.Range("M80").Value = 282
$ws.Range("N80").Value = -3291.8182

$ws.Range("H83").Value = 1141.2
$ws.Range("I83").Value = 716
$ws.Range("J83").Value = 1295.8182
$ws.Range("K83").Value = 3580
$ws.Range("L83").Value = 6479.090999999999
$ws.Range("M83").Value = 1412
$ws.Range("N83").Value = -16463.091

$ws.Range("H98").Value = 99999
$ws.Range("J98").Value = 99999
$ws.Range("L98").Value = 99999
$ws.Range("N98").Value = -105989

$ws.Range("H134").Value = 2441159.2
$ws.Range("I134").Value = 2018.5555
$ws.Range("K134").Value = 6055.666499999999
$ws.Range("M134").Value = -3520.666499999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 79.95238000000001
$ws.Range("I7").Value = 226
$ws.Range("J7").Value = 55.61111
$ws.Range("K7").Value = 226
$ws.Range("L7").Value = 55.61111
$ws.Range("M7").Value = -113
$ws.Range("N7").Value = -281.61111

$ws.Range("H31").Value = 21279274
$ws.Range("I31").Value = 30305218
$ws.Range("J31").Value = 3831.2856
$ws.Range("K31").Value = 30305218
$ws.Range("L31").Value = 3831.2856
$ws.Range("M31").Value = -30304923
$ws.Range("N31").Value = -4421.2856

$ws.Range("H34").Value = 21279274
$ws.Range("I34").Value = 30305218
$ws.Range("J34").Value = 3831.2856
$ws.Range("K34").Value = 30305218
$ws.Range("L34").Value = 3831.2856
$ws.Range("M34").Value = -30305016
$ws.Range("N34").Value = -4235.2856

$ws.Range("H99").Value = 13051.682
$ws.Range("I99").Value = 6567.077
$ws.Range("J99").Value = 22418.334
$ws.Range("K99").Value = 6567.077
$ws.Range("L99").Value = 22418.334
$ws.Range("M99").Value = -5069.077
$ws.Range("N99").Value = -25414.334

$ws.Range("H103").Value = 34963.777
$ws.Range("J103").Value = 60749.5
$ws.Range("L103").Value = 60749.5
$ws.Range("N103").Value = -63093.5

$ws.Range("H107").Value = 1281.591
$ws.Range("I107").Value = 423.64706
$ws.Range("K107").Value = 423.64706
$ws.Range("M107").Value = 1496.35294

$ws.Range("H126").Value = 13051.682
$ws.Range("I126").Value = 6567.077
$ws.Range("J126").Value = 22418.334
$ws.Range("K126").Value = 19701.231
$ws.Range("L126").Value = 67255.00199999999
$ws.Range("M126").Value = -17231.231
$ws.Range("N126").Value = -72195.00199999999

$ws.Range("H134").Value = 1351.02
$ws.Range("I134").Value = 1079.4103
$ws.Range("K134").Value = 3238.2309
$ws.Range("M134").Value = -703.2309

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 9206.4
$ws.Range("I17").Value = 500
$ws.Range("J17").Value = 11383
$ws.Range("K17").Value = 1500
$ws.Range("L17").Value = 34149
$ws.Range("M17").Value = -1331
$ws.Range("N17").Value = -34487

$ws.Range("H60").Value = 6938.909
$ws.Range("I60").Value = 385
$ws.Range("K60").Value = 1155
$ws.Range("M60").Value = -904

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 38328.332
$ws.Range("I5").Value = 38328.332
$ws.Range("K5").Value = 38328.332
$ws.Range("M5").Value = -38216.332

$ws.Range("H18").Value = 49999
$ws.Range("I18").Value = 49999
$ws.Range("K18").Value = 49999
$ws.Range("M18").Value = -49706

$ws.Range("H52").Value = 18979.8
$ws.Range("I52").Value = 10000
$ws.Range("J52").Value = 24966.334
$ws.Range("K52").Value = 10000
$ws.Range("L52").Value = 24966.334
$ws.Range("M52").Value = -9741
$ws.Range("N52").Value = -25484.334

$ws.Range("H57").Value = 24624.5
$ws.Range("J57").Value = 38999.11
$ws.Range("L57").Value = 38999.11
$ws.Range("N57").Value = -40639.11

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H23").Value = 4982.25
$ws.Range("I23").Value = 4982.25
$ws.Range("K23").Value = 4982.25
$ws.Range("M23").Value = -4752.25

$ws.Range("H40").Value = 5461.421
$ws.Range("I40").Value = 4523.769
$ws.Range("K40").Value = 4523.769
$ws.Range("M40").Value = -4387.769

$ws.Range("H61").Value = 2679.8823
$ws.Range("I61").Value = 2397.5925
$ws.Range("K61").Value = 2397.5925
$ws.Range("M61").Value = -2195.5925

$ws.Range("H74").Value = 96666
$ws.Range("I74").Value = 94999.5
$ws.Range("K74").Value = 94999.5
$ws.Range("M74").Value = -94001.5

$ws.Range("H77").Value = 96666
$ws.Range("I77").Value = 94999.5
$ws.Range("K77").Value = 284998.5
$ws.Range("M77").Value = -280006.5

$ws.Range("H113").Value = 2679.8823
$ws.Range("I113").Value = 2397.5925
$ws.Range("K113").Value = 2397.5925
$ws.Range("M113").Value = -227.5925000000002

$ws.Range("H122").Value = 3999.7437
$ws.Range("I122").Value = 3499.8333
$ws.Range("J122").Value = 9998.666999999999
$ws.Range("K122").Value = 10499.4999
$ws.Range("L122").Value = 29996.001
$ws.Range("M122").Value = -8049.499899999999
$ws.Range("N122").Value = -34896.001

$ws.Range("H132").Value = 4048.756
$ws.Range("I132").Value = 2156.5454
$ws.Range("K132").Value = 6469.6362
$ws.Range("M132").Value = -3939.6362

$ws.Range("H136").Value = 4377.1577
$ws.Range("I136").Value = 1973.1666
$ws.Range("K136").Value = 5919.4998
$ws.Range("M136").Value = -3369.4998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H93").Value = 99999
$ws.Range("J93").Value = 99999
$ws.Range("L93").Value = 99999
$ws.Range("N93").Value = -104991

$ws.Range("H113").Value = 365.5238
$ws.Range("I113").Value = 182.3077
$ws.Range("K113").Value = 546.9231
$ws.Range("M113").Value = 1623.0769

$ws.Range("H122").Value = 3079.2144
$ws.Range("I122").Value = 2430.5
$ws.Range("K122").Value = 7291.5
$ws.Range("M122").Value = -4841.5

$ws.Range("H126").Value = 6541.3335
$ws.Range("I126").Value = 6859
$ws.Range("J126").Value = 4000
$ws.Range("K126").Value = 20577
$ws.Range("L126").Value = 12000
$ws.Range("M126").Value = -18107
$ws.Range("N126").Value = -16940

$ws.Range("H132").Value = 215496.06
$ws.Range("I132").Value = 2716.628
$ws.Range("K132").Value = 8149.884
$ws.Range("M132").Value = -5619.884

